$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 260.9091
$ws.Cells.Item(19, 9).Value = 399.6
$ws.Cells.Item(19, 10).Value = 145.33333
$ws.Cells.Item(19, 11).Value = 399.6
$ws.Cells.Item(19, 12).Value = 145.33333
$ws.Cells.Item(19, 13).Value = -224.6
$ws.Cells.Item(19, 14).Value = -495.33333
$ws.Cells.Item(33, 8).Value = 6055.8823
$ws.Cells.Item(33, 9).Value = 50.285713
$ws.Cells.Item(33, 10).Value = 10259.8
$ws.Cells.Item(33, 11).Value = 50.285713
$ws.Cells.Item(33, 12).Value = 10259.8
$ws.Cells.Item(33, 13).Value = 178.714287
$ws.Cells.Item(33, 14).Value = -10717.8
$ws.Cells.Item(40, 8).Value = 1670.909
$ws.Cells.Item(40, 9).Value = 1573.7333
$ws.Cells.Item(40, 10).Value = 1879.1428
$ws.Cells.Item(40, 11).Value = 1573.7333
$ws.Cells.Item(40, 12).Value = 1879.1428
$ws.Cells.Item(40, 13).Value = -1398.7333
$ws.Cells.Item(40, 14).Value = -2229.1428
$ws.Cells.Item(68, 8).Value = 15000
$ws.Cells.Item(68, 10).Value = 15000
$ws.Cells.Item(68, 12).Value = 15000
$ws.Cells.Item(68, 14).Value = -16498
$ws.Cells.Item(71, 8).Value = 15000
$ws.Cells.Item(71, 10).Value = 15000
$ws.Cells.Item(71, 12).Value = 45000
$ws.Cells.Item(71, 14).Value = -52488
$ws.Cells.Item(80, 8).Value = 1962.9762
$ws.Cells.Item(80, 9).Value = 488.2
$ws.Cells.Item(80, 10).Value = 4131.7646
$ws.Cells.Item(80, 11).Value = 1464.6
$ws.Cells.Item(80, 12).Value = 12395.2938
$ws.Cells.Item(80, 13).Value = -466.5999999999999
$ws.Cells.Item(80, 14).Value = -14391.2938
$ws.Cells.Item(83, 8).Value = 1962.9762
$ws.Cells.Item(83, 9).Value = 488.2
$ws.Cells.Item(83, 10).Value = 4131.7646
$ws.Cells.Item(83, 11).Value = 4393.8
$ws.Cells.Item(83, 12).Value = 37185.88140000001
$ws.Cells.Item(83, 13).Value = 598.1999999999998
$ws.Cells.Item(83, 14).Value = -47169.88140000001
$ws.Cells.Item(109, 8).Value = 30000
$ws.Cells.Item(109, 10).Value = 30000
$ws.Cells.Item(109, 12).Value = 30000
$ws.Cells.Item(109, 14).Value = -32774
$ws.Cells.Item(132, 8).Value = 1741.5834
$ws.Cells.Item(132, 9).Value = 1644.3529
$ws.Cells.Item(132, 11).Value = 4933.0587
$ws.Cells.Item(132, 13).Value = -2403.0587

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 3723.0557
$ws.Cells.Item(45, 9).Value = 4080.8
$ws.Cells.Item(45, 10).Value = 3275.875
$ws.Cells.Item(45, 11).Value = 4080.8
$ws.Cells.Item(45, 12).Value = 3275.875
$ws.Cells.Item(45, 13).Value = -3703.8
$ws.Cells.Item(45, 14).Value = -4029.875
$ws.Cells.Item(61, 8).Value = 2432.1052
$ws.Cells.Item(61, 9).Value = 2435
$ws.Cells.Item(61, 10).Value = 2421.25
$ws.Cells.Item(61, 11).Value = 2435
$ws.Cells.Item(61, 12).Value = 2421.25
$ws.Cells.Item(61, 13).Value = -2223
$ws.Cells.Item(61, 14).Value = -2845.25
$ws.Cells.Item(122, 8).Value = 3516.6155
$ws.Cells.Item(122, 9).Value = 3669.0952
$ws.Cells.Item(122, 10).Value = 2876.2
$ws.Cells.Item(122, 11).Value = 11007.2856
$ws.Cells.Item(122, 12).Value = 8628.599999999999
$ws.Cells.Item(122, 13).Value = -8557.285600000001
$ws.Cells.Item(122, 14).Value = -13528.6
$ws.Cells.Item(123, 8).Value = 32000
$ws.Cells.Item(123, 10).Value = 32000
$ws.Cells.Item(123, 12).Value = 32000
$ws.Cells.Item(123, 14).Value = -41800
$ws.Cells.Item(136, 8).Value = 2432.1052
$ws.Cells.Item(136, 9).Value = 2435
$ws.Cells.Item(136, 10).Value = 2421.25
$ws.Cells.Item(136, 11).Value = 7305
$ws.Cells.Item(136, 12).Value = 7263.75
$ws.Cells.Item(136, 13).Value = -4755
$ws.Cells.Item(136, 14).Value = -12363.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 1653.775
$ws.Cells.Item(105, 9).Value = 1603.3334
$ws.Cells.Item(105, 11).Value = 1603.3334
$ws.Cells.Item(105, 13).Value = 143.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 1962.7736
$ws.Cells.Item(99, 9).Value = 1932.6888
$ws.Cells.Item(99, 10).Value = 2132
$ws.Cells.Item(99, 11).Value = 1932.6888
$ws.Cells.Item(99, 12).Value = 2132
$ws.Cells.Item(99, 13).Value = -434.6887999999999
$ws.Cells.Item(99, 14).Value = -5128
$ws.Cells.Item(126, 8).Value = 1962.7736
$ws.Cells.Item(126, 9).Value = 1932.6888
$ws.Cells.Item(126, 10).Value = 2132
$ws.Cells.Item(126, 11).Value = 5798.0664
$ws.Cells.Item(126, 12).Value = 6396
$ws.Cells.Item(126, 13).Value = -3328.0664
$ws.Cells.Item(126, 14).Value = -11336

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 147
$ws.Cells.Item(12, 9).Value = 35.333332
$ws.Cells.Item(12, 10).Value = 172.76923
$ws.Cells.Item(12, 11).Value = 105.999996
$ws.Cells.Item(12, 12).Value = 518.30769
$ws.Cells.Item(12, 13).Value = 67.000004
$ws.Cells.Item(12, 14).Value = -864.30769

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 2412.9333
$ws.Cells.Item(122, 9).Value = 2357.8333
$ws.Cells.Item(122, 10).Value = 2633.3333
$ws.Cells.Item(122, 11).Value = 7073.499899999999
$ws.Cells.Item(122, 12).Value = 7899.999899999999
$ws.Cells.Item(122, 13).Value = -4623.499899999999
$ws.Cells.Item(122, 14).Value = -12799.9999
$ws.Cells.Item(126, 8).Value = 27781270
$ws.Cells.Item(126, 9).Value = 4190
$ws.Cells.Item(126, 10).Value = 166666670
$ws.Cells.Item(126, 11).Value = 12570
$ws.Cells.Item(126, 12).Value = 500000010
$ws.Cells.Item(126, 13).Value = -10100
$ws.Cells.Item(126, 14).Value = -500004950

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 5804.591
$ws.Cells.Item(46, 9).Value = 1366.75
$ws.Cells.Item(46, 10).Value = 11130
$ws.Cells.Item(46, 11).Value = 1366.75
$ws.Cells.Item(46, 12).Value = 11130
$ws.Cells.Item(46, 13).Value = -1178.75
$ws.Cells.Item(46, 14).Value = -11506
$ws.Cells.Item(68, 8).Value = 2230.4348
$ws.Cells.Item(68, 9).Value = 2221.4285
$ws.Cells.Item(68, 10).Value = 2244.4443
$ws.Cells.Item(68, 11).Value = 2221.4285
$ws.Cells.Item(68, 12).Value = 2244.4443
$ws.Cells.Item(68, 13).Value = -1472.4285
$ws.Cells.Item(68, 14).Value = -3742.4443
$ws.Cells.Item(71, 8).Value = 2230.4348
$ws.Cells.Item(71, 9).Value = 2221.4285
$ws.Cells.Item(71, 10).Value = 2244.4443
$ws.Cells.Item(71, 11).Value = 11107.1425
$ws.Cells.Item(71, 12).Value = 11222.2215
$ws.Cells.Item(71, 13).Value = -7363.1425
$ws.Cells.Item(71, 14).Value = -18710.2215
$ws.Cells.Item(100, 8).Value = 3444.8696
$ws.Cells.Item(100, 9).Value = 3172.2727
$ws.Cells.Item(100, 11).Value = 3172.2727
$ws.Cells.Item(100, 13).Value = -2631.2727
$ws.Cells.Item(122, 8).Value = 2698.1428
$ws.Cells.Item(122, 9).Value = 2483.3333
$ws.Cells.Item(122, 10).Value = 3987
$ws.Cells.Item(122, 11).Value = 7449.999899999999
$ws.Cells.Item(122, 12).Value = 11961
$ws.Cells.Item(122, 13).Value = -4999.999899999999
$ws.Cells.Item(122, 14).Value = -16861
$ws.Cells.Item(136, 8).Value = 3099.4546
$ws.Cells.Item(136, 9).Value = 3011.75
$ws.Cells.Item(136, 10).Value = 3333.3333
$ws.Cells.Item(136, 11).Value = 9035.25
$ws.Cells.Item(136, 12).Value = 9999.999899999999
$ws.Cells.Item(136, 13).Value = -6485.25
$ws.Cells.Item(136, 14).Value = -15099.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 1454.3846
$ws.Cells.Item(122, 9).Value = 1450.5834
$ws.Cells.Item(122, 10).Value = 1500
$ws.Cells.Item(122, 11).Value = 4351.7502
$ws.Cells.Item(122, 12).Value = 4500
$ws.Cells.Item(122, 13).Value = -1901.7502
$ws.Cells.Item(122, 14).Value = -9400
